# Auto-generated edit script: update market-data snapshot values
# across the 8 crafting-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR),
# mirroring a scheduled market-data refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 8999.875
$ws.Range("I86").Value = 7800
$ws.Range("K86").Value = 7800
$ws.Range("M86").Value = -6677

$ws.Range("H89").Value = 8999.875
$ws.Range("I89").Value = 7800
$ws.Range("K89").Value = 39000
$ws.Range("M89").Value = -33384

$ws.Range("H138").Value = 1484.1333
$ws.Range("I138").Value = 1294
$ws.Range("J138").Value = 1864.4
$ws.Range("K138").Value = 3882
$ws.Range("L138").Value = 5593.200000000001
$ws.Range("M138").Value = 1258
$ws.Range("N138").Value = -15873.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1665.3334
$ws.Range("I61").Value = 1698.8
$ws.Range("K61").Value = 1698.8
$ws.Range("M61").Value = -1486.8

$ws.Range("H132").Value = 1186.4166
$ws.Range("I132").Value = 1186.4166
$ws.Range("K132").Value = 3559.2498
$ws.Range("M132").Value = -1029.2498

$ws.Range("H136").Value = 1665.3334
$ws.Range("I136").Value = 1698.8
$ws.Range("K136").Value = 5096.4
$ws.Range("M136").Value = -2546.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1974.5

$ws.Range("H67").Value = 1974.5

$ws.Range("H75").Value = 1375
$ws.Range("I75").Value = 1375
$ws.Range("K75").Value = 1375
$ws.Range("M75").Value = -439

$ws.Range("H78").Value = 1375
$ws.Range("I78").Value = 1375
$ws.Range("K78").Value = 4125
$ws.Range("M78").Value = 555

$ws.Range("H94").Value = 2149.1428
$ws.Range("I94").Value = 2149.1428
$ws.Range("K94").Value = 2149.1428
$ws.Range("M94").Value = -1698.1428

$ws.Range("H107").Value = 991.3333
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1441.35
$ws.Range("I7").Value = 789.1875
$ws.Range("K7").Value = 789.1875
$ws.Range("M7").Value = -676.1875

$ws.Range("H31").Value = 3102.5
$ws.Range("I31").Value = 3102.5
$ws.Range("K31").Value = 3102.5
$ws.Range("M31").Value = -2807.5

$ws.Range("H34").Value = 3102.5
$ws.Range("I34").Value = 3102.5
$ws.Range("K34").Value = 3102.5
$ws.Range("M34").Value = -2900.5

$ws.Range("H92").Value = 41247.8
$ws.Range("J92").Value = 41247.8
$ws.Range("L92").Value = 41247.8
$ws.Range("N92").Value = -46239.8

$ws.Range("H107").Value = 783.2778
$ws.Range("I107").Value = 868.8
$ws.Range("J107").Value = 355.66666
$ws.Range("K107").Value = 868.8
$ws.Range("L107").Value = 355.66666
$ws.Range("M107").Value = 1051.2
$ws.Range("N107").Value = -4195.66666

$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 450
$ws.Range("I21").Value = 450
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 1350
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -1177
$ws.Range("N21").ClearContents()

$ws.Range("H97").Value = 587.2857
$ws.Range("I97").Value = 676.75
$ws.Range("K97").Value = 2030.25
$ws.Range("M97").Value = -1534.25

$ws.Range("H107").Value = 1652
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()

$ws.Range("H129").Value = 1283
$ws.Range("I129").Value = 596.2
$ws.Range("J129").Value = 3000
$ws.Range("K129").Value = 1788.6
$ws.Range("L129").Value = 9000
$ws.Range("M129").Value = 3211.4
$ws.Range("N129").Value = -19000

$ws.Range("H131").Value = 985.8333

$ws.Range("H139").Value = 3943
$ws.Range("I139").Value = 3943
$ws.Range("K139").Value = 11829
$ws.Range("M139").Value = -6689

$ws.Range("H140").Value = 1499.5
$ws.Range("I140").Value = 1500
$ws.Range("J140").Value = 1499
$ws.Range("K140").Value = 4500
$ws.Range("L140").Value = 4497
$ws.Range("M140").Value = 680
$ws.Range("N140").Value = -14857

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws.Range("H97").Value = 3609
$ws.Range("I97").Value = 3577.4
$ws.Range("J97").Value = 3672.2
$ws.Range("K97").Value = 3577.4
$ws.Range("L97").Value = 3672.2
$ws.Range("M97").Value = -3081.4
$ws.Range("N97").Value = -4664.2

$ws.Range("H132").Value = 4769.857
$ws.Range("I132").Value = 4769.857
$ws.Range("K132").Value = 14309.571
$ws.Range("M132").Value = -11779.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2290.2727
$ws.Range("I22").Value = 1100
$ws.Range("J22").Value = 2970.4285
$ws.Range("K22").Value = 1100
$ws.Range("L22").Value = 2970.4285
$ws.Range("M22").Value = -805
$ws.Range("N22").Value = -3560.4285

$ws.Range("H27").Value = 2290.2727
$ws.Range("I27").Value = 1100
$ws.Range("J27").Value = 2970.4285
$ws.Range("K27").Value = 1100
$ws.Range("L27").Value = 2970.4285
$ws.Range("M27").Value = -993
$ws.Range("N27").Value = -3184.4285

$ws.Range("H93").Value = 1999.75
$ws.Range("I93").Value = 1999.75
$ws.Range("K93").Value = 1999.75
$ws.Range("M93").Value = -751.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 5198.5
$ws.Range("J4").Value = 3598
$ws.Range("L4").Value = 3598
$ws.Range("N4").Value = -3824

$ws.Range("H5").Value = 5000
$ws.Range("J5").Value = 5000
$ws.Range("L5").Value = 5000
$ws.Range("N5").Value = -5224

$ws.Range("H96").Value = 1990.909
$ws.Range("I96").Value = 1987.5
$ws.Range("J96").Value = 2000
$ws.Range("K96").Value = 1987.5
$ws.Range("L96").Value = 2000
$ws.Range("M96").Value = -614.5
$ws.Range("N96").Value = -4746

$ws.Range("H122").Value = 1540.9565
$ws.Range("I122").Value = 1618
$ws.Range("K122").Value = 4854
$ws.Range("M122").Value = -2404

$ws.Range("H132").Value = 3583.3333
$ws.Range("I132").Value = 4500
$ws.Range("J132").Value = 3400
$ws.Range("K132").Value = 13500
$ws.Range("L132").Value = 10200
$ws.Range("M132").Value = -10970
$ws.Range("N132").Value = -15260
